$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI")

# Rename header labels (row 1) — parentheses/spaces swapped for underscores
$ws.Range("B1").Value = "What_Action_Items"
$ws.Range("D1").Value = "Who_Responsible"
$ws.Range("E1").Value = "When_Due_date"

# Update the "VoC project" action text (row 8, column F) with revised deadlines
$ws.Range("F8").Value = "1. Need discussion with Yamazaki-san and other party to solve user's question by end of December`n2. Send one VoC survey after providing answers to users by end of January 2022"
$ws.Rows.Item(8).RowHeight = 43.5

# Update the "Knowledge base chatbot project" action text (row 12, column F) with revised deadline
$ws.Range("F12").Value = "1. Have to discuss further with IT representative by mid of January`n"

# Move the active selection to E1
$ws.Range("E1").Select()
